$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).Insert()
$ws.Range("A3").Value = "Match ID"
$ws.Range("A4:A19").Value = 10
$ws.Range("A20").Value = 10
$ws.Range("A3:A19").Font.Bold = $true
$ws.Range("A3:A19").Select()
$excel.ActiveWindow.ScrollRow = 3
